# Replace the body content with the new résumé-style summary text.
# Using Range.InsertXML lets us author the exact run/paragraph/proofErr
# structure (incl. a manual line break and two trailing empty paragraphs)
# in one shot, the way Word's clipboard/XML-paste pipeline would.
$d = $word.ActiveDocument
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Jeg er en </w:t></w:r><w:r><w:t>ambitiøs</w:t></w:r><w:r><w:t xml:space="preserve"> kandidat – studerende der er motiveret for at lærer mere inden for programmerings verden. På nuværende tidspunkt er jeg på mit 3. semester på kandidaten. Jeg søger et job hvor jeg kan anvende min opnået viden til at programmere løsninger. Jeg er målrettet og byder udfordringer velkommen. </w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Jeg har erfaring indenfor diverse programmeringssprog som f.eks. Java, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Python</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, C#, R og </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Javascript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Herudover har jeg designet forskellige frontends til applikationer så jeg er bekendt med HTML, CSS (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bootstrap</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) og adskillige JavaScript </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>frameworks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Jeg har udviklet løsninger ved brug af pattern</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> såsom MVC og har specielt viden indenfor </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>den objekt</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> orienteret tilgang. På mit arbejde og </w:t></w:r><w:r><w:t xml:space="preserve">i </w:t></w:r><w:r><w:t>nogle få projekter har jeg udviklet modeller både til business intelligence og machine (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>deep</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>learning</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Jeg håber på at anvende ovenstående til at udfylde en organisations opgaver og mål i fremtiden. Glæder mig til at høre fra jer!</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
